# Fill previously-blank numeric placeholder cells with a "-" inline string,
# matching the formatting already used elsewhere in the sheet for
# "not applicable" values (right-aligned text with thin border, style index 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18 and 23: columns D through O were empty numeric placeholders.
$ranges = @("D18:O18", "D23:O23")

# Rows 26 and 33: columns D through F already show "-"; only G through O
# were empty numeric placeholders.
$ranges += "G26:O26"
$ranges += "G33:O33"

foreach ($addr in $ranges) {
    $rng = $ws.Range($addr)
    $rng.Value = "-"
    $rng.HorizontalAlignment = -4152   # xlRight
    $rng.Borders.LineStyle = 1         # xlContinuous (thin border, matches existing "-" cells)
}
